# WP-08: finalize stochastic pricing integration and recalibration
#
# Recalibrates the senior-tranche "sculpted" DSCR target from 1.45x to
# 1.20x: the descriptive label on the Intro sheet, the frozen/"sculpted"
# DSCR schedule on the Debt Service sheet (target column N and the
# pre-solved sculpted values in column K), the mirrored schedule on the
# Sculpted_Waterfall sheet (column Z) and the underlying DSCR_Obj
# parameter on the "Params (waterfall only)" sheet (B10).
#
# The K/Z columns hold previously goal-seeked ("sculpted") values that
# are stored as plain numbers rather than live formulas (the workbook
# does not have iterative calculation enabled, so those circular
# sculpting formulas were solved once and pasted as values). Switching
# Calculation to manual for the duration of the edits -- then restoring
# it to automatic -- lets us rebase those cached numbers to the new
# 1.20x target without Excel trying to re-solve the circular "sculpted"
# waterfall chain (which would cascade zeros through the dependent
# amortization columns).

$wb = $excel.ActiveWorkbook
$wb.Application.Calculation = -4135  # xlCalculationManual

# --- Intro sheet: update the descriptive label ------------------------
$introWs = $wb.Worksheets.Item("Intro")
$introWs.Range("B21").Value = "“Sculpted” para mantener DSCR ≈ 1.20 en senior"

# --- Debt Service sheet: recalibrate target (N) and sculpted (K) cols -
$dsWs = $wb.Worksheets.Item("Debt Service")

# Column N ("Objetivo") target DSCR, rows 3-17, all set to 1.20
$dsWs.Range("N3").Value = 1.2
$dsWs.Range("N4").Value = 1.2
$dsWs.Range("N5").Value = 1.2
$dsWs.Range("N6").Value = 1.2
$dsWs.Range("N7").Value = 1.2
$dsWs.Range("N8").Value = 1.2
$dsWs.Range("N9").Value = 1.2
$dsWs.Range("N10").Value = 1.2
$dsWs.Range("N11").Value = 1.2
$dsWs.Range("N12").Value = 1.2
$dsWs.Range("N13").Value = 1.2
$dsWs.Range("N14").Value = 1.2
$dsWs.Range("N15").Value = 1.2
$dsWs.Range("N16").Value = 1.2
$dsWs.Range("N17").Value = 1.2

# Column K sculpted DSCR values - only the entries that were previously
# ~1.45 (rebased to ~1.20); rows that were already off-target (K7, K10)
# keep their original untouched values.
$dsWs.Range("K8").Value = 1.2
$dsWs.Range("K9").Value = 1.200000000032805
$dsWs.Range("K11").Value = 1.2000000001719371
$dsWs.Range("K12").Value = 1.2
$dsWs.Range("K13").Value = 1.2

# --- Sculpted_Waterfall sheet: recalibrate sculpted DSCR (Z) column ---
$swWs = $wb.Worksheets.Item("Sculpted_Waterfall")
$swWs.Range("Z7").Value = 1.2
$swWs.Range("Z8").Value = 1.200000000032805
$swWs.Range("Z10").Value = 1.2000000001719371
$swWs.Range("Z11").Value = 1.2
$swWs.Range("Z12").Value = 1.2

# --- Params (waterfall only) sheet: recalibrate target DSCR parameter -
$paramsWs = $wb.Worksheets.Item("Params (waterfall only)")
$paramsWs.Range("B10").Value = 1.2

$wb.Application.Calculation = -4105  # xlCalculationAutomatic
